# "Add files via upload" - new vocabulary rows appended to the word list on
# Sheet1, directly below the existing "majority / 大多数" entry (row 181).
# Rows 182-197 were previously blank (the sheet already had data again at
# row 198), so the nine new word pairs just fill that gap - no row
# insertion/shifting is required.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# length / 长度  (row 182) - typed Chinese first, then the English term.
$ws.Range("B182").Value = "长度"
$ws.Range("A183").Value = "interrupt"
$ws.Range("B183").Value = "打断"
$ws.Range("A184").Value = "enable"
$ws.Range("B184").Value = "使可能"
$ws.Range("A185").Value = "selection"
$ws.Range("B185").Value = "选择"
$ws.Range("A182").Value = "length"
$ws.Range("A186").Value = "trigger"
$ws.Range("B186").Value = "触发"
$ws.Range("A187").Value = "modem"
$ws.Range("B187").Value = "调制解调器"
$ws.Range("A188").Value = "overrun"
$ws.Range("B188").Value = "超支"
$ws.Range("A189").Value = "utility"
$ws.Range("B189").Value = "效用"
$ws.Range("A190").Value = "facility"
$ws.Range("B190").Value = "设施"

# Leave the cursor where the author ended up after typing the new rows.
$ws.Range("C185").Select() | Out-Null
